$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 41 (shifts existing row 41 "im" and everything
# below it down by one), then populate the new row with the new dictionary
# entry "idh2015" / "Indice de desarrollo humano, 2015".
$ws.Rows.Item(41).Insert()

$ws.Cells.Item(41, 1).Value = "idh2015"
$ws.Cells.Item(41, 2).Value = "Índice de desarrollo humano, 2015"

# Update the view: scroll so row 22 is the top-left row, and select B42
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("B42").Select()
